$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("device-logistics")

# Insert a new row at row 10 for the new "team project" device entry,
# pushing the existing patch/charger logistics rows down by one.
$ws.Rows("10").Insert()

$ws.Range("A10").Value = "phone_samsung"
$ws.Range("B10").Value = "J3"
$ws.Range("C10").Value = "ANDROID7"
$ws.Range("D10").Value = "Pratik"
$ws.Range("E10").Value = "Unable to log-in"

# Update the logistics for the most recently lost device
# (charger_vivalnk / C700136), now shifted down to row 20.
$ws.Range("D20").Value = "LOST"

# E20 previously had no cell at all, so give it the same style as its
# row-mates before writing the status note.
$ws.Range("D20").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("E20").Value = "MRN: 1303033"

# Leave the selection where the author left it after the edit.
$ws.Range("A11").Select()
